$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing 6-column data block (A1:F28) before restructuring.
$orig = $ws.Range("A1:F28").Value()

$paraplyNames = @{
  "MAVIT" = "MAVIT Helsevitenskap"
  "MASY"  = "MASY Spesialsykepleie"
  "MSFH"  = "MSFH Sosialfag (H)"
  "MSFD"  = "MSFD Sosialfag (D)"
}

# Build the new 7-column grid: A,B unchanged; new C = Paraplynamn; old C,D,E,F shift to D,E,F,G.
$new = New-Object 'object[,]' 28,7
for ($r = 0; $r -lt 28; $r++) {
  $codeA = $orig[$r+1,1]
  $codeB = $orig[$r+1,2]

  $new[$r,0] = $codeA
  $new[$r,1] = $codeB
  if ($r -eq 0) {
    $new[$r,2] = "Paraplynamn"
  } else {
    $new[$r,2] = $paraplyNames[$codeA]
  }
  $new[$r,3] = $orig[$r+1,3]
  $new[$r,4] = $orig[$r+1,4]
  $new[$r,5] = $orig[$r+1,5]
  $new[$r,6] = $orig[$r+1,6]
}

$ws.Range("A1:G28").Value = $new

Write-Host "grid done"
